$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet has an unused "Ghi chu" / notes column (column I) that needs to be
# removed entirely. Removing it shifts every column from J onward one to the
# left (J->I, K->J, L->K, ... T->S) and drops the now-unused "Ghi chu" shared
# string.
#
# Comments anchored in row 1 of the columns to the right of the deleted
# column need to keep their text but move to the new (shifted) cell
# reference, so we capture their text up front, remove the stale comment
# objects, delete the column, then re-create the comments on the
# correctly-shifted cells.

$commentCells = @("J1", "L1", "P1", "Q1", "R1", "S1", "T1")

$savedComments = @{}
foreach ($addr in $commentCells) {
    $rng = $ws.Range($addr)
    if ($null -ne $rng.Comment) {
        $savedComments[$addr] = $rng.Comment.Text()
        $rng.Comment.Delete()
    }
}

# Delete the whole column I (the "Ghi chu" column) and shift everything
# after it to the left.
$ws.Columns.Item(9).Delete()

# Map each original comment anchor to its new (post-deletion) anchor.
$moveMap = @{
    "J1" = "I1"
    "L1" = "K1"
    "P1" = "O1"
    "Q1" = "P1"
    "R1" = "Q1"
    "S1" = "R1"
    "T1" = "S1"
}

foreach ($oldAddr in $commentCells) {
    if ($savedComments.ContainsKey($oldAddr)) {
        $newAddr = $moveMap[$oldAddr]
        $text = $savedComments[$oldAddr]
        $newRng = $ws.Range($newAddr)
        $newComment = $newRng.AddComment($text)
        try {
            $newComment.Shape.TextFrame.Characters(1, 6).Font.Bold = $true
        } catch {
        }
    }
}

# Restore the selection that Excel stores for the sheet view.
$ws.Activate()
$ws.Range("F9").Select()
